# Imports.xlsx — "ajout groupes dans fichier d'imports"
# Adds a new "INFOS-GROUPES" worksheet (id / libelle / parcours) at the
# end of the workbook and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet at the end of the tab strip ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "INFOS-GROUPES"

# --- Header row ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "libelle"
$ws.Range("C1").Value = "parcours"

# --- Data rows (write order chosen so the shared-string table is built
#     up in the same order as the source workbook) ---
$ws.Range("B2").Value = "A"
$ws.Range("C2").Value = "RAPP_6"
$ws.Range("A2").Value = "in_s5_A_2026_2027"

$ws.Range("A3").Value = "TestImportsMultiples"
$ws.Range("B3").Value = "A"
$ws.Range("C3").Value = "DACS"

# --- Cosmetics: fit column A to its new content, put the cursor on A5 ---
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Range("A5").Select() | Out-Null

# --- Make the new sheet the active tab (also clears tabSelected on
#     whichever sheet was previously active) ---
$ws.Activate() | Out-Null
